# Journal de bord - add a new day (D column: 2021-10-06) of entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New text entries --------------------------------------------------
# Written in this exact order so the new shared-string table entries are
# appended in the same order as in the target workbook.
$ws.Range("D2").Value = "Ajout de liens a la BDD et rerévision"
$ws.Range("D3").Value = "commencement a intégrer un formlaire connexion & inscription"
$ws.Range("D4").Value = "API rest en reflexion & non fonctionnel"
$ws.Range("D7").Value = "que mettre dans l'API rest"
$ws.Range("D8").Value = "des liens dans la BDD difficile a comprendre"
$ws.Range("D9").Value = "idcompte champ de tableau invalide (liens mal fait)"
$ws.Range("D5").Value = "page inscription/connexion fonctionnel"
$ws.Range("D6").Value = "page de planning fonctionnel"

# --- New date header (D1), formatted like the other date cells --------
$ws.Range("D1").Value = 44475
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column widths --------------------------------------------------
# ColumnWidth is quantized to whole pixels by the engine (same as real
# Excel), so these are the closest achievable values to the stored
# widths of 58.140625 / 43.85546875 characters.
$ws.Columns.Item(3).ColumnWidth = 57.25
$ws.Columns.Item(4).ColumnWidth = 42.92

# --- Selection restore (active cell ends on D5) -------------------------
$ws.Range("D5").Select() | Out-Null
